# Update workbook/sheet for the "through 11-10" data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet and its title text (2021-11-09 -> 2021-11-10)
$ws.Name = "Through 2021-11-10"

# Update the "November (through 11-09)" label to "November (through 11-10)"
$ws.Range("A12").Value = "November (through 11-10)"

# Update November row (row 12) values for columns C:H (2016-2021)
$ws.Range("C12").Value = 23
$ws.Range("D12").Value = 37
$ws.Range("E12").Value = 25
$ws.Range("F12").Value = 17
$ws.Range("G12").Value = 63
$ws.Range("H12").Value = 68

# Update Total row (row 13) values for columns C:H (2016-2021)
$ws.Range("C13").Value = 509
$ws.Range("D13").Value = 747
$ws.Range("E13").Value = 640
$ws.Range("F13").Value = 499
$ws.Range("G13").Value = 1120
$ws.Range("H13").Value = 1512
